# "suppression de compétence obsolète"
# Remove the obsolete skill codes A3.3.3, A3.3.4 and A3.2.2 from the
# "IV Listes des activités" block, leaving a double-tab gap where each
# code used to sit (mirrors the author's edit, which left the ";"
# separators/tabs in place but dropped the stale reference), and fix the
# typo A3.3.2 -> A3.2.2 on the last line.

$d = $word.ActiveDocument

# NOTE: these are typed as [string] (not [char]) on purpose -- in this
# PowerShell dialect "[char] + [char]" performs *numeric* addition, not
# string concatenation, which would silently corrupt every literal built
# below out of two or more adjacent NBSP/TAB characters.
$NBSP = [string][char]160
$TAB  = [string][char]9

# --- Paragraph "A1.1.1 ; A1.2.3 ; A1.4.1 ; A3.3.3 ; A4.1.4 ; A5.1.1" ---
# Drop "A3.3.3" (keep its leading separator + a double tab before "; A4.1.4").
$old1 = $NBSP + "; A3.3.3" + $NBSP + "; A4.1.4"
$new1 = $NBSP + $TAB + $TAB + "; A4.1.4"
$found = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Paragraph "A1.1.3 ; A1.2.5 ; A2.3.2 ; A3.3.4 ; A4.1.5 ; A5.2.3" ---
# This line has a _GoBack bookmark sitting between " ; A" and "1.2.5", so the
# replacement is done in two pieces (one on either side of the bookmark) to
# avoid disturbing it, exactly as the source XML is split.
$old2a = "A1.1.3" + $NBSP + "; A"
$new2a = "A1.1.3" + $NBSP + "; A1.2.5" + $NBSP + "; A2.3.2" + $NBSP + $NBSP + ";" + $TAB + $TAB
$found = $d.Content.Find.Execute($old2a, $false, $false, $false, $false, $false, $true, 1, $false, $new2a, 2)

$old2b = "1.2.5" + $NBSP + "; A2.3.2" + $NBSP + "; A3.3.4" + $NBSP + "; A4.1.5"
$new2b = $NBSP + "A4.1.5"
$found = $d.Content.Find.Execute($old2b, $false, $false, $false, $false, $false, $true, 1, $false, $new2b, 2)

# --- Paragraph "A1.2.1 ; A1.3.2 ; A3.2.2 ; A4.1.3 ; A4.1.7 ; A5.2.4" ---
# Drop "A3.2.2" (keep the separator/tabs before "A4.1.3").
$old3 = $NBSP + "; A3.2.2" + $NBSP + "; A4.1.3"
$new3 = $NBSP + $NBSP + ";" + $NBSP + $TAB + $TAB + "A4.1.3"
$found = $d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# --- Paragraph "A1.2.2 ; A1.3.4 ; A3.3.2 ; A4.1.2 ; A4.2.1" ---
# Fix the typo "A3.3.2" -> "A3.2.2".
$old4 = "A3.3.2"
$new4 = "A3.2.2"
$found = $d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
